$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.933
$ws.Range("B21").Value = 9.379000000000001
$ws.Range("B23").Value = 7.886
$ws.Range("D24").Value = -6.968999999999999
$ws.Range("B25").Value = 5.392999999999999
$ws.Range("D28").Value = -8.235999999999999
$ws.Range("D36").Value = -7.540999999999999
$ws.Range("D45").Value = -7.425000000000002
$ws.Range("D48").Value = -7.265000000000001
$ws.Range("D49").Value = -8.322999999999999
$ws.Range("D52").Value = -8.080000000000002
$ws.Range("B53").Value = 5.706
$ws.Range("D53").Value = -8.343999999999999
$ws.Range("D54").Value = -8.168000000000001
$ws.Range("B57").Value = 5.034000000000001
$ws.Range("B59").Value = 4.790000000000001
$ws.Range("B69").Value = 5.724
$ws.Range("D70").Value = -7.026999999999999
$ws.Range("B79").Value = 5.889
$ws.Range("B83").Value = 5.473999999999999
$ws.Range("D86").Value = -8.252000000000001
$ws.Range("D87").Value = -8.233999999999998
$ws.Range("B93").Value = 5.369
$ws.Range("D101").Value = -7.806999999999999
